$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.632.32"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.564.79"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "210.73"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +6.71%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "24.81"
$ws.Range("E8").Value = "  +5.54%  "
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "0.0899"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.787.23"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "1.591.87"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "28.665.84"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "61.65"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "227.10"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "7.30"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "0.0₃0684"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "151.92"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "0.108"
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").Value = "14.79"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "6.26"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "1.403.61"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").Value = "2.72"
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "0.517"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").Value = "0.0459"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "63.91"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  -8.96%  "
$ws.Range("D49").Value = "84.68"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "42.16"
$ws.Range("E50").Value = "  +3.65%  "
$ws.Range("E51").Value = "  -0.66%  "
